$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 34 (Athens Kallithea vs Giouchtas -> Ionikos vs Tilikratis L.) ---
# Row 34
$ws.Range("F34").Value = 'Ionikos'
$ws.Range("H34").Value = 'Tilikratis L.'
$ws.Range("L34").Value = 1.22
$ws.Range("M34").Value = '21/10/2023 12:52'
$ws.Range("P34").Value = 5.42
$ws.Range("Q34").Value = '21/10/2023 13:02'
$ws.Range("T34").Value = 14.23
$ws.Range("U34").Value = '21/10/2023 12:52'
$ws.Range("V34").Value = 'https://www.betexplorer.com/football/greece/super-league-2/ionikos-tilikratis-lefkada/hry7rFPN/'

# Row 35
$ws.Range("F35").Value = 'PAOK B'
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 'Aiolikos'
$ws.Range("J35").Value = 1.57
$ws.Range("L35").Value = 1.75
$ws.Range("M35").Value = '21/10/2023 14:43'
$ws.Range("N35").Value = 3.91
$ws.Range("P35").Value = 3.66
$ws.Range("Q35").Value = '21/10/2023 14:43'
$ws.Range("R35").Value = 5.06
$ws.Range("T35").Value = 4.44
$ws.Range("U35").Value = '21/10/2023 14:43'
$ws.Range("V35").Value = 'https://www.betexplorer.com/football/greece/super-league-2/paok-aiolikos-fc/fyiX48y2/'

# Row 36
$ws.Range("F36").Value = 'Athens Kallithea'
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = 'Giouchtas'
$ws.Range("J36").Value = 1.3
$ws.Range("L36").Value = 1.23
$ws.Range("M36").Value = '21/10/2023 14:17'
$ws.Range("N36").Value = 4.55
$ws.Range("P36").Value = 5.44
$ws.Range("Q36").Value = '21/10/2023 14:18'
$ws.Range("R36").Value = 10.53
$ws.Range("T36").Value = 13.52
$ws.Range("U36").Value = '21/10/2023 14:17'
$ws.Range("V36").Value = 'https://www.betexplorer.com/football/greece/super-league-2/athens-kallithea-giouchtas/n5zBsZvU/'

# Row 39
$ws.Range("F39").Value = 'Panachaiki'
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 'PAE Egaleo'
$ws.Range("I39").Value = 1
$ws.Range("J39").Value = 3.03
$ws.Range("L39").Value = 3.1
$ws.Range("M39").Value = '22/10/2023 14:59'
$ws.Range("N39").Value = 2.94
$ws.Range("P39").Value = 3.06
$ws.Range("Q39").Value = '22/10/2023 14:57'
$ws.Range("R39").Value = 2.3
$ws.Range("T39").Value = 2.37
$ws.Range("U39").Value = '22/10/2023 14:59'
$ws.Range("V39").Value = 'https://www.betexplorer.com/football/greece/super-league-2/panachaiki-pae-egaleo/QV5L1xQo/'

# Row 41
$ws.Range("F41").Value = 'Apollon Pontou'
$ws.Range("H41").Value = 'Niki Volos'
$ws.Range("I41").Value = 3
$ws.Range("J41").Value = 4.17
$ws.Range("K41").Value = '21/10/2023 02:13'
$ws.Range("L41").Value = 9.99
$ws.Range("M41").Value = '22/10/2023 14:56'
$ws.Range("N41").Value = 3.18
$ws.Range("O41").Value = '21/10/2023 02:13'
$ws.Range("P41").Value = 4.44
$ws.Range("Q41").Value = '22/10/2023 14:56'
$ws.Range("R41").Value = 1.78
$ws.Range("S41").Value = '21/10/2023 02:13'
$ws.Range("T41").Value = 1.34
$ws.Range("U41").Value = '22/10/2023 14:56'
$ws.Range("V41").Value = 'https://www.betexplorer.com/football/greece/super-league-2/apollon-pontou-niki-volos/UTgL7Aiq/'

# Row 42
$ws.Range("F42").Value = 'Kozani FC'
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = 'Karditsa'
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 2.75
$ws.Range("K42").Value = '22/10/2023 12:12'
$ws.Range("L42").Value = 2.58
$ws.Range("M42").Value = '22/10/2023 14:49'
$ws.Range("N42").Value = 2.74
$ws.Range("O42").Value = '22/10/2023 12:12'
$ws.Range("P42").Value = 2.74
$ws.Range("Q42").Value = '22/10/2023 13:54'
$ws.Range("R42").Value = 2.92
$ws.Range("S42").Value = '22/10/2023 12:12'
$ws.Range("T42").Value = 3.14
$ws.Range("U42").Value = '22/10/2023 14:49'
$ws.Range("V42").Value = 'https://www.betexplorer.com/football/greece/super-league-2/kozani-fc-karditsa/OnHhCjqS/'

# --- Add new row 45 ---
$ws.Range("A44:V44").Copy() | Out-Null
$ws.Range("A45:V45").PasteSpecial(-4122) | Out-Null
$ws.Range("A45").Value = 44
$ws.Range("B45").Value = 'greece'
$ws.Range("C45").Value = 'super-league-2'
$ws.Range("D45").Value = '2023-2024'
$ws.Range("E45").Value = 45224.625
$ws.Range("F45").Value = 'Iraklis 1908'
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 'PAOK B'
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 1.61
$ws.Range("K45").Value = '14/10/2023 02:12'
$ws.Range("L45").Value = 1.69
$ws.Range("M45").Value = '25/10/2023 14:51'
$ws.Range("N45").Value = 3.58
$ws.Range("O45").Value = '14/10/2023 02:12'
$ws.Range("P45").Value = 3.65
$ws.Range("Q45").Value = '25/10/2023 14:51'
$ws.Range("R45").Value = 4.69
$ws.Range("S45").Value = '14/10/2023 02:12'
$ws.Range("T45").Value = 4.82
$ws.Range("U45").Value = '25/10/2023 14:52'
$ws.Range("V45").Value = 'https://www.betexplorer.com/football/greece/super-league-2/iraklis-fc-paok/GMDpECEF/'
